$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Row 15 (field #7 "FinalFg") gains a remark/comment explaining the update rule.
$ws.Range("G15").Value = "同一戶號下最近申請的案件編號其下之關係人會更新成Y"

# Row 17 (CreateDate) and row 19 (LastUpdate) data type corrected from DATE to TIMESTAMP.
$ws.Range("D17").Value = "TIMESTAMP"
$ws.Range("D19").Value = "TIMESTAMP"

# Reflect the updated selection/scroll position left by the editor.
$ws.Activate() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("D19").Select() | Out-Null
